# The presentation's slide master currently uses the "Integral" theme
# (colour scheme: dk1=000000, lt1=FFFFFF, dk2=455F51, lt2=E3DED1,
# accent1=99CB38, accent2=63A537, accent3=E6D024, accent4=CC9700,
# accent5=4EB3CF, accent6=378DA6, hlink=6B9F25, folHlink=B26B02).
#
# The target edit swaps the deck's active design back to the default
# Office theme colour palette:
# dk1=000000, lt1=FFFFFF, dk2=44546A, lt2=E7E6E6, accent1=5B9BD5,
# accent2=ED7D31, accent3=A5A5A5, accent4=FFC000, accent5=4472C4,
# accent6=70AD47, hlink=0563C1, folHlink=954F72.
#
# Font scheme (Arial/Arial) and format scheme are already identical
# between the two themes, so only the twelve theme colours need to
# change. PowerPoint's theme colour slots are addressed in a fixed
# 1..12 order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$themeColors = $master.Theme.ThemeColorScheme

$newColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

for ($i = 1; $i -le $newColors.Count; $i++) {
    $hex = $newColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $rgb = $r -bor ($g -shl 8) -bor ($b -shl 16)
    $themeColors.Colors($i).RGB = $rgb
}
